# The two Pearson Edexcel logo pictures (one embedded in each footer)
# are renamed from "image1.png" to "image2.png", and the BTEC logo
# picture (embedded in the first-page header) is renamed from
# "image2.jpg" to "image1.jpg".
#
# Headers/footers aren't reliably reachable/writable through
# Section.Headers/Section.Footers in this host, so walk every story
# range instead and identify each inline picture by its (stable)
# AlternativeText / description.

$d = $word.ActiveDocument

foreach ($story in $d.StoryRanges) {
    $shapes = $story.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $ishape = $shapes.Item($i)
        $descr = $ishape.AlternativeText

        $newName = $null
        if ($descr -like "*PearsonLogo.png") {
            $newName = "image2.png"
        } elseif ($descr -eq "BTec_Logo-Orange") {
            $newName = "image1.jpg"
        }

        if ($newName -ne $null) {
            # Renaming straight on the InlineShape is unreliable for
            # footer stories in this host, but the equivalent
            # ShapeRange accessor renames it correctly everywhere.
            $ishape.Range.ShapeRange.Item(1).Name = $newName
        }
    }
}

Write-Output "Renamed logo inline shapes."
